# Apply updated cryptocurrency price/volume data to Sheet1
# (values refreshed by the periodic GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.314.13'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '1.664.74'
$ws.Range('D4').Value = '''1.010'
$ws.Range('E4').Value = '  +0.92%  '
$ws.Range('D5').Value = '''218.88'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').Value = '''0.5342'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('D8').Value = '''0.2659'
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('D9').Value = '''0.06400'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').Value = '''0.07835'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '''4.566'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').Value = '1.664.52'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = '1.892.56'
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').Value = '''0.5525'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '0.0₅8230'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').Value = '''65.78'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').Value = '''4.695'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '''10.25'
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('D22').Value = '''6.040'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '''1.012'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').Value = '''146.30'
$ws.Range('E24').Value = '  +3.03%  '
$ws.Range('D25').Value = '''0.1231'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').Value = '''7.195'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').Value = '''16.10'
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('D28').Value = '''1.485'
$ws.Range('E28').Value = '  +3.97%  '
$ws.Range('D29').Value = '''0.05841'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').Value = '''1.284'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('D31').Value = '''3.625'
$ws.Range('E31').Value = '  +2.54%  '
$ws.Range('D32').Value = '''3.280'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').Value = '''1.617'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').Value = '''0.9663'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('E35').Value = '  +1.64%  '
$ws.Range('D36').Value = '''2.418'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').Value = '''0.5810'
$ws.Range('E37').Value = '  +1.72%  '
$ws.Range('D38').Value = '''0.01605'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').Value = '''0.8694'
$ws.Range('E39').Value = '  +2.73%  '
$ws.Range('D40').Value = '''5.873'
$ws.Range('E40').Value = '  +1.59%  '
$ws.Range('D41').Value = '1.052.50'
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '''105.13'
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '''1.010'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('D44').Value = '1.803.99'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').Value = '''57.95'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').Value = '''1.014'
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈104'
$ws.Range('E47').Value = '  -6.64%  '
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('D49').Value = '''7.997'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').Value = '''1.411'
$ws.Range('E51').Value = '  -3.94%  '

Write-Host "Applied 98 cell updates."
